# Update the multiplication problems in the single table of the document.
# Each data cell's text is replaced in-place using the Tables/Cell object
# model so that overlapping old/new values (e.g. "340×3=" is both a
# replacement target and a replacement source) cannot cause incorrect
# cascading matches.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Map of (row, column) -> new value, in row-major order matching the
# five data rows of the table (rows 1, 5, 10, 15, 20).
$updates = @(
    @{Row=1;  Col=1; New="480×7="},
    @{Row=1;  Col=2; New="221×8="},
    @{Row=1;  Col=3; New="590×5="},
    @{Row=1;  Col=4; New="706×3="},
    @{Row=1;  Col=5; New="459×8="},

    @{Row=5;  Col=1; New="395×6="},
    @{Row=5;  Col=2; New="726×2="},
    @{Row=5;  Col=3; New="924×9="},
    @{Row=5;  Col=4; New="340×3="},
    @{Row=5;  Col=5; New="817×7="},

    @{Row=10; Col=1; New="332×5="},
    @{Row=10; Col=2; New="869×9="},
    @{Row=10; Col=3; New="399×4="},
    @{Row=10; Col=4; New="836×4="},
    @{Row=10; Col=5; New="612×6="},

    @{Row=15; Col=1; New="629×9="},
    @{Row=15; Col=2; New="152×4="},
    @{Row=15; Col=3; New="951×5="},
    @{Row=15; Col=4; New="879×7="},
    @{Row=15; Col=5; New="646×5="},

    @{Row=20; Col=1; New="338×4="},
    @{Row=20; Col=2; New="202×7="},
    @{Row=20; Col=3; New="321×6="},
    @{Row=20; Col=4; New="562×5="},
    @{Row=20; Col=5; New="847×2="}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $r = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, preserving the run's formatting.
    $r.End = $r.End - 1
    $r.Text = $u.New
}
